$ws = $excel.ActiveWorkbook.ActiveSheet
$ws.Range("D2").Value = '29.133.12'
$ws.Range("D3").Value = '1.833.30'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.27'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6575'
$ws.Range("E6").Value = '  -1.19%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07398'
$ws.Range("E8").Value = '  +0.56%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2925'
$ws.Range("E9").Value = '  -1.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.84'
$ws.Range("E10").Value = '  +0.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07743'
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("D12").Value = '1.843.92'
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.987'
$ws.Range("E13").Value = '  -0.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6664'
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '83.17'
$ws.Range("E15").Value = '  -3.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.118'
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008606'
$ws.Range("E17").Value = '  +4.37%  '
$ws.Range("D18").Value = '29.143.44'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = '2.082.08'
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '226.38'
$ws.Range("E20").Value = '  -1.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.44'
$ws.Range("E21").Value = '  -0.62%  '
$ws.Range("E22").Value = '  +0.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.100'
$ws.Range("E23").Value = '  -2.80%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '161.08'
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1402'
$ws.Range("E26").Value = '  -1.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.573'
$ws.Range("E27").Value = '  -1.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.97'
$ws.Range("E28").Value = '  -0.27%  '
$ws.Range("E29").Value = '  +0.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.100'
$ws.Range("E30").Value = '  -3.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.040'
$ws.Range("E31").Value = '  -1.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.189'
$ws.Range("E32").Value = '  -1.18%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05280'
$ws.Range("E33").Value = '  -0.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.862'
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7372'
$ws.Range("E35").Value = '  -1.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.142'
$ws.Range("E36").Value = '  +1.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.657'
$ws.Range("E37").Value = '  -0.82%  '
$ws.Range("D38").Value = '1.300.56'
$ws.Range("E38").Value = '  -1.24%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01787'
$ws.Range("E39").Value = '  -1.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.738'
$ws.Range("E40").Value = '  +0.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9155'
$ws.Range("E41").Value = '  -1.09%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.005'
$ws.Range("E42").Value = '  +0.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9997'
$ws.Range("E43").Value = '  +0.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.08161'
$ws.Range("E44").Value = '  +9.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.26'
$ws.Range("E45").Value = '  -1.22%  '
$ws.Range("D46").Value = '1.985.62'
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5132'
$ws.Range("E47").Value = '  -0.72%  '
$ws.Range("E48").Value = '  -0.63%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '63.59'
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.746'
$ws.Range("E50").Value = '  -0.95%  '
$ws.Range("E51").Value = '  -1.41%  '
